$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.238222002983093
$ws.Range("B1").Value = 2.492285251617432
$ws.Range("C1").Value = 3.919533252716064
$ws.Range("D1").Value = 2.775002002716064
$ws.Range("E1").Value = 1.084582090377808
